{"js": "const replacements = [\n  [\"514\u00d74=2056\", \"422\u00d78=3376\"],\n  [\"566\u00d76=3396\", \"848\u00d76=5088\"],\n  [\"164\u00d76=984\", \"536\u00d78=4288\"],\n  [\"778\u00d79=7002\", \"777\u00d73=2331\"],\n  [\"830\u00d72=1660\", \"989\u00d78=7912\"],\n  [\"275\u00d79=2475\", \"978\u00d76=5868\"],\n  [\"413\u00d72=826\", \"161\u00d78=1288\"],\n  [\"479\u00d77=3353\", \"639\u00d74=2556\"],\n  [\"721\u00d79=6489\", \"105\u00d78=840\"],\n  [\"275\u00d76=1650\", \"280\u00d78=2240\"],\n  [\"533\u00d77=3731\", \"953\u00d73=2859\"],\n  [\"785\u00d74=3140\", \"699\u00d78=5592\"],\n  [\"333\u00d72=666\", \"523\u00d79=4707\"],\n  [\"348\u00d74=1392\", \"878\u00d79=7902\"],\n  [\"319\u00d79=2871\", \"200\u00d79=1800\"],\n  [\"277\u00d77=1939\", \"506\u00d75=2530\"],\n  [\"430\u00d78=3440\", \"966\u00d74=3864\"],\n  [\"480\u00d75=2400\", \"584\u00d75=2920\"],\n  [\"217\u00d76=1302\", \"753\u00d79=6777\"],\n  [\"664\u00d79=5976\", \"993\u00d74=3972\"],\n  [\"389\u00d77=2723\", \"964\u00d73=2892\"],\n  [\"547\u00d78=4376\", \"534\u00d75=2670\"],\n  [\"378\u00d79=3402\", \"413\u00d74=1652\"],\n  [\"616\u00d73=1848\", \"396\u00d74=1584\"],\n  [\"699\u00d73=2097\", \"837\u00d79=7533\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$wdReplaceNone = 0\n$wdReplaceOne = 1\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"514\u00d74=2056\", \"422\u00d78=3376\"),\n  @(\"566\u00d76=3396\", \"848\u00d76=5088\"),\n  @(\"164\u00d76=984\", \"536\u00d78=4288\"),\n  @(\"778\u00d79=7002\", \"777\u00d73=2331\"),\n  @(\"830\u00d72=1660\", \"989\u00d78=7912\"),\n  @(\"275\u00d79=2475\", \"978\u00d76=5868\"),\n  @(\"413\u00d72=826\", \"161\u00d78=1288\"),\n  @(\"479\u00d77=3353\", \"639\u00d74=2556\"),\n  @(\"721\u00d79=6489\", \"105\u00d78=840\"),\n  @(\"275\u00d76=1650\", \"280\u00d78=2240\"),\n  @(\"533\u00d77=3731\", \"953\u00d73=2859\"),\n  @(\"785\u00d74=3140\", \"699\u00d78=5592\"),\n  @(\"333\u00d72=666\", \"523\u00d79=4707\"),\n  @(\"348\u00d74=1392\", \"878\u00d79=7902\"),\n  @(\"319\u00d79=2871\", \"200\u00d79=1800\"),\n  @(\"277\u00d77=1939\", \"506\u00d75=2530\"),\n  @(\"430\u00d78=3440\", \"966\u00d74=3864\"),\n  @(\"480\u00d75=2400\", \"584\u00d75=2920\"),\n  @(\"217\u00d76=1302\", \"753\u00d79=6777\"),\n  @(\"664\u00d79=5976\", \"993\u00d74=3972\"),\n  @(\"389\u00d77=2723\", \"964\u00d73=2892\"),\n  @(\"547\u00d78=4376\", \"534\u00d75=2670\"),\n  @(\"378\u00d79=3402\", \"413\u00d74=1652\"),\n  @(\"616\u00d73=1848\", \"396\u00d74=1584\"),\n  @(\"699\u00d73=2097\", \"837\u00d79=7533\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n  if (-not $found) {\n    Write-Output \"WARNING: not found -> $oldText\"\n  }\n}\n"}
